# Insert a new data row at row 184 (pushes the existing rows 184-238 down to
# 185-239) and populate it with the new weekly record for "Ajo" / "Chino" /
# "Primera" from Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 184..238 down by one row.
$ws.Rows(184).Insert()

# Populate the newly inserted row 184 with the new record.
$ws.Cells.Item(184, 1).Value = 9
$ws.Cells.Item(184, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(184, 3).Value = "Metropolitana"
$ws.Cells.Item(184, 4).Value = 44711
$ws.Cells.Item(184, 5).Value = 13
$ws.Cells.Item(184, 6).Value = 100112003
$ws.Cells.Item(184, 7).Value = "Ajo"
$ws.Cells.Item(184, 8).Value = "Chino"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 610
$ws.Cells.Item(184, 11).Value = 18000
$ws.Cells.Item(184, 12).Value = 18500
$ws.Cells.Item(184, 13).Value = 18250
$ws.Cells.Item(184, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(184, 15).Value = "China"
$ws.Cells.Item(184, 16).Value = 1825
$ws.Cells.Item(184, 17).Value = 10
$ws.Cells.Item(184, 18).Value = "Hortaliza"
